# Update the workbook to reflect the latest Twitter keyword used for the
# Google Alerts / BoardReader / Twitter scraping run.
#
# Content change: cell B10 ("Twitter Keyword:" row) is updated from the
# previous hashtag to the new one, "#istandwithraeesah".
#
# We also move the active selection to B10, matching where the user was
# last working in the sheet when they made this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B10").Value = "#istandwithraeesah"

$ws.Range("B10").Select()
